$wb = $excel.ActiveWorkbook

# --- Reorder sheets: move "Roof" to sit right after "Facade" ---
$roof = $wb.Worksheets.Item("Roof")
$roof.Move($null, $wb.Worksheets.Item("Facade"))

# --- Rename sheets: the old "Roof" data becomes "Ground"; the old
#     "Ground(Off,Ret)"/"Ground(Res)" sheets become "Roof(Off,Ret)"/"Roof(Res)" ---
$wb.Worksheets.Item("Roof").Name = "Ground"
$wb.Worksheets.Item("Ground(Off,Ret)").Name = "Roof(Off,Ret)"
$wb.Worksheets.Item("Ground(Res)").Name = "Roof(Res)"

# --- Re-enter B32 on Facade as an explicit (non-shared) formula ---
$facade = $wb.Worksheets.Item("Facade")
$facade.Range("B32").Formula = "=1/(0.360091+A32/0.037)"

# --- Reset every sheet's view to B1 (clears any stored topLeftCell),
#     activating "Roof(Res)" last so it ends up the selected tab ---
$tabOrder = @("Facade", "Ground", "Roof(Off,Ret)", "Window", "Roof(Res)")
foreach ($name in $tabOrder) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate()
    $ws.Range("B1").Select()
}
